$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 104
    $ws.Range("F5").Value = 2676
    $ws.Range("F6").Value = 253
}
